# Updating collaboration.R with 2022 data
#
# The underlying rows in the citations sheet were re-ordered: the row that
# used to sit at row 3 now belongs at row 4 (and vice versa), and likewise
# for rows 7 and 8. This script swaps the full contents of those row pairs
# in place, cell by cell, preserving each cell's original data type
# (text / boolean / number) and avoiding any lasting style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel (real COM semantics, mirrored here) infers a cell's type from the
# text you hand it - a numeric- or date-looking string silently becomes a
# number/date serial. The source workbook stores these columns as literal
# text, so force "Text" number-format while assigning, then restore the
# cell style so no stray formatting is left behind.
function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

function Swap-TextCells {
    param($Ws, $Addr1, $Addr2)
    $r1 = $Ws.Range($Addr1)
    $r2 = $Ws.Range($Addr2)
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    Set-TextValue $r1 $v2
    Set-TextValue $r2 $v1
}

function Swap-BoolCells {
    param($Ws, $Addr1, $Addr2)
    $r1 = $Ws.Range($Addr1)
    $r2 = $Ws.Range($Addr2)
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

# ---- Rows 3 <-> 4 (full swap; every populated column matches on both sides) ----
$textCols34 = @("A","B","C","E","F","G","H","J","K","L","O","P","Q","R","U","V","AC","AE")
foreach ($col in $textCols34) {
    Swap-TextCells $ws "$($col)3" "$($col)4"
}
Swap-BoolCells $ws "W3" "W4"

# ---- Rows 7 <-> 8 ----
# Same idea, except row 7 originally carried an abstract in column E while
# row 8 had none - after the swap that abstract (and its presence/absence)
# moves along with the rest of the row 7 content into row 8.
$textCols78 = @("A","B","C","F","G","H","J","K","O","P","Q","AC","AE")
foreach ($col in $textCols78) {
    Swap-TextCells $ws "$($col)7" "$($col)8"
}

$e7 = $ws.Range("E7").Value2
$ws.Range("E7").ClearContents()
Set-TextValue $ws.Range("E8") $e7
